$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ConfirmationID" column (G) used for email confirmation tracking.

# G1: header cell, copy the header style/formatting from an existing header
# cell (A1) and then set its text.
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "ConfirmationID"

# G2: materialize an (empty) cell in the new column, matching the plain
# (unstyled) formatting used by the other data cells in row 2.
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# G3: confirmation id value for the second row of data.
$ws.Range("G3").Value = "94a8008e-c7ed-4642-9526-df6ed7af3261"

$excel.CutCopyMode = 0
